$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for 3 new articles at the top of the data (row 2) ---
$ws.Rows("2:4").Insert()
# Insert() copies the formatting of the row above (the bold/centered header row);
# reset the newly inserted rows back to plain/default formatting like the other data rows.
$ws.Rows("2:4").ClearFormats()

# --- Step 2: drop the 3 oldest articles that fall off the bottom of the list ---
# (after the insert above, the old rows 9,10,11 are now at 12,13,14)
$ws.Rows("12:14").Delete()

# --- Step 3: fill the newly inserted rows 2-4 with the new articles ---
$newRows = @(
    @{A="2026-01-28"; B="2026-01-27"; C="OpenAI"; D="PVH reimagines the future of fashion with OpenAI"; E="PVH는 OpenAI를 통해 패션의 미래를 재구상합니다."; F="https://openai.com/index/pvh-future-of-fashion"},
    @{A="2026-01-28"; B="2026-01-27"; C="OpenAI"; D="Powering tax donations with AI powered personalized recommendations"; E="AI 기반 맞춤형 추천으로 세금 기부 지원"; F="https://openai.com/index/trustbank"},
    @{A="2026-01-28"; B="2026-01-27"; C="OpenAI"; D="Introducing Prism"; E="프리즘 소개"; F="https://openai.com/index/introducing-prism"}
)

$r = 2
foreach ($item in $newRows) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $r = $r + 1
}

# --- Step 4: update the "collection date" (column A) for the carried-over rows (now 5-11) ---
for ($row = 5; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = "2026-01-28"
}

# --- Step 5: row 7 ("Scaling PostgreSQL...") also received a freshly re-worded Korean translation ---
$ws.Cells.Item(7, 5).Value = "8억 명의 ChatGPT 사용자를 지원하기 위해 PostgreSQL 확장"
